# Applies the "Atualizacao de bases das ligas" data refresh to the
# "Poland IV Liga" sheet: a handful of match rows get corrected
# results/odds, three pairs of rows trade places (same date/pool,
# ids resorted), and a few HomeTeam/AwayTeam cells are re-pointed so
# the displayed club names stay correct after that resort.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Poland IV Liga")

# Row 6
$ws.Range("B6").Value = 6750018
$ws.Range("F6").Value = 'Swit Starozreby'
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 'D'
$ws.Range("J6").Value = 1.571
$ws.Range("L6").Value = 4.333
$ws.Range("M6").Value = 1.4
$ws.Range("N6").Value = 4.5
$ws.Range("O6").Value = 6
$ws.Range("P6").Value = -1.25
$ws.Range("Q6").Value = 1.8
$ws.Range("R6").Value = 2
$ws.Range("T6").Value = 1.8
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = -1
$ws.Range("W6").Value = 3.5
$ws.Range("Y6").Value = -1
$ws.Range("Z6").Value = 1
$ws.Range("AB6").Value = 1

# Row 7
$ws.Range("B7").Value = 6746871
$ws.Range("F7").Value = 'Chemik Bydgoszcz'
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 'H'
$ws.Range("J7").Value = 3.25
$ws.Range("L7").Value = 1.8
$ws.Range("M7").Value = 3.5
$ws.Range("N7").Value = 4
$ws.Range("O7").Value = 1.727
$ws.Range("P7").Value = 0.75
$ws.Range("Q7").Value = 1.825
$ws.Range("R7").Value = 1.975
$ws.Range("T7").Value = 1.9
$ws.Range("U7").Value = 1.9
$ws.Range("V7").Value = 2.5
$ws.Range("W7").Value = -1
$ws.Range("Y7").Value = 0.825
$ws.Range("Z7").Value = -1
$ws.Range("AB7").Value = 0.8999999999999999

# Row 16
$ws.Range("F16").Value = 'Tarnovia Tarnowo Podgorne'

# Row 32
$ws.Range("E32").Value = 'Hutnik Warsaw'

# Row 43
$ws.Range("B43").Value = 7068602
$ws.Range("J43").Value = 2.1
$ws.Range("K43").Value = 4.2
$ws.Range("L43").Value = 2.5
$ws.Range("M43").Value = 2.1
$ws.Range("N43").Value = 4.2
$ws.Range("O43").Value = 2.55
$ws.Range("P43").Value = -0.25
$ws.Range("Q43").Value = 1.95
$ws.Range("R43").Value = 1.85
$ws.Range("S43").Value = 3
$ws.Range("T43").Value = 1.775
$ws.Range("U43").Value = 2.025
$ws.Range("W43").Value = 3.2
$ws.Range("Y43").Value = -0.5
$ws.Range("Z43").Value = 0.425
$ws.Range("AA43").Value = 0.7749999999999999

# Row 44
$ws.Range("B44").Value = 7068603
$ws.Range("J44").Value = 3.75
$ws.Range("K44").Value = 4.333
$ws.Range("L44").Value = 1.615
$ws.Range("M44").Value = 3.75
$ws.Range("N44").Value = 4.333
$ws.Range("O44").Value = 1.615
$ws.Range("P44").Value = 0.75
$ws.Range("Q44").Value = 2
$ws.Range("R44").Value = 1.8
$ws.Range("S44").Value = 3.25
$ws.Range("T44").Value = 2.025
$ws.Range("U44").Value = 1.775
$ws.Range("W44").Value = 3.333
$ws.Range("Y44").Value = 1
$ws.Range("Z44").Value = -1
$ws.Range("AA44").Value = 1.025

# Row 45
$ws.Range("B45").Value = 7068599
$ws.Range("E45").Value = 'LKS Jawiszowice'
$ws.Range("F45").Value = 'MKS TrzebiniaSiersza'
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 1.5
$ws.Range("K45").Value = 4.5
$ws.Range("L45").Value = 4.333
$ws.Range("M45").Value = 1.5
$ws.Range("N45").Value = 4.5
$ws.Range("O45").Value = 4.333
$ws.Range("P45").Value = -1
$ws.Range("Q45").Value = 1.8
$ws.Range("R45").Value = 2
$ws.Range("V45").Value = 0.5
$ws.Range("Y45").Value = 0
$ws.Range("Z45").Value = 0
$ws.Range("AA45").Value = -1
$ws.Range("AB45").Value = 0.95

# Row 46
$ws.Range("B46").Value = 7068596
$ws.Range("E46").Value = 'MGKS Moto Jelcz Olawa'
$ws.Range("F46").Value = 'Iskra Ksieginice'
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 1
$ws.Range("I46").Value = 'H'
$ws.Range("J46").Value = 2
$ws.Range("L46").Value = 2.75
$ws.Range("M46").Value = 2
$ws.Range("O46").Value = 2.8
$ws.Range("P46").Value = -0.25
$ws.Range("Q46").Value = 1.825
$ws.Range("R46").Value = 1.975
$ws.Range("S46").Value = 3.25
$ws.Range("V46").Value = 1
$ws.Range("X46").Value = -1
$ws.Range("Y46").Value = 0.825
$ws.Range("Z46").Value = -1
$ws.Range("AA46").Value = 0.8500000000000001
$ws.Range("AB46").Value = -1

# Row 47
$ws.Range("B47").Value = 7068595
$ws.Range("E47").Value = 'Prochowiczanka Prochowice'
$ws.Range("F47").Value = 'Piast Nowa Ruda'
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 2
$ws.Range("I47").Value = 'A'
$ws.Range("J47").Value = 1.95
$ws.Range("K47").Value = 4
$ws.Range("L47").Value = 2.875
$ws.Range("M47").Value = 1.95
$ws.Range("N47").Value = 4
$ws.Range("O47").Value = 2.9
$ws.Range("P47").Value = -0.5
$ws.Range("Q47").Value = 2
$ws.Range("R47").Value = 1.8
$ws.Range("S47").Value = 3.5
$ws.Range("V47").Value = -1
$ws.Range("X47").Value = 1.9
$ws.Range("Y47").Value = -1
$ws.Range("Z47").Value = 0.8

# Row 72
$ws.Range("B72").Value = 7183411
$ws.Range("E72").Value = 'Piast Gliwice II'
$ws.Range("F72").Value = 'Zaglebie Sosnowiec II'
$ws.Range("H72").Value = 3
$ws.Range("J72").Value = 1.909
$ws.Range("K72").Value = 4
$ws.Range("L72").Value = 2.9
$ws.Range("M72").Value = 1.909
$ws.Range("N72").Value = 4
$ws.Range("O72").Value = 2.9
$ws.Range("P72").Value = -0.5
$ws.Range("Q72").Value = 1.95
$ws.Range("R72").Value = 1.85
$ws.Range("X72").Value = 1.9
$ws.Range("Z72").Value = 0.8500000000000001
$ws.Range("AA72").Value = 0.8500000000000001
$ws.Range("AB72").Value = -1

# Row 73
$ws.Range("B73").Value = 7183410
$ws.Range("E73").Value = 'Arka Gdynia II'
$ws.Range("F73").Value = 'Grom Nowy Staw'
$ws.Range("G73").Value = 1
$ws.Range("J73").Value = 2.6
$ws.Range("K73").Value = 3.6
$ws.Range("L73").Value = 2.2
$ws.Range("M73").Value = 2.6
$ws.Range("N73").Value = 3.6
$ws.Range("O73").Value = 2.2
$ws.Range("P73").Value = 0
$ws.Range("Q73").Value = 2.05
$ws.Range("R73").Value = 1.75
$ws.Range("X73").Value = 1.2
$ws.Range("Z73").Value = 0.75

# Row 74
$ws.Range("B74").Value = 7183407
$ws.Range("E74").Value = 'Prochowiczanka Prochowice'
$ws.Range("F74").Value = 'Miedz Legnica II'
$ws.Range("G74").Value = 0
$ws.Range("H74").Value = 2
$ws.Range("J74").Value = 2.25
$ws.Range("K74").Value = 3.75
$ws.Range("L74").Value = 2.5
$ws.Range("M74").Value = 4.5
$ws.Range("N74").Value = 4.2
$ws.Range("O74").Value = 1.533
$ws.Range("P74").Value = 1
$ws.Range("Q74").Value = 1.975
$ws.Range("R74").Value = 1.825
$ws.Range("X74").Value = 0.5329999999999999
$ws.Range("Z74").Value = 0.825
$ws.Range("AA74").Value = -1
$ws.Range("AB74").Value = 0.95

# Row 77
$ws.Range("E77").Value = 'Hutnik Warsaw'

# Row 102
$ws.Range("F102").Value = 'Tarnovia Tarnowo Podgorne'

# Row 109
$ws.Range("E109").Value = 'Korona Piaski'

Write-Output "Applied 188 cell updates across 15 rows"
